# Append the 2025-06-10 USDA families resazurin size measurements
# (25 samples: A1-A5, B1-B5, C1-C5, D1-D5, E1-E5) as rows 127-151.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(127, 1).Value = 20250610
$ws.Cells.Item(127, 2).Value = "A1"
$ws.Cells.Item(127, 3).Value = 13.942
$ws.Cells.Item(127, 4).Value = 127.848

$ws.Cells.Item(128, 1).Value = 20250610
$ws.Cells.Item(128, 2).Value = "A2"
$ws.Cells.Item(128, 3).Value = 18.236000000000001
$ws.Cells.Item(128, 4).Value = 198.29900000000001

$ws.Cells.Item(129, 1).Value = 20250610
$ws.Cells.Item(129, 2).Value = "A3"
$ws.Cells.Item(129, 3).Value = 13.949
$ws.Cells.Item(129, 4).Value = 125.464

$ws.Cells.Item(130, 1).Value = 20250610
$ws.Cells.Item(130, 2).Value = "A4"
$ws.Cells.Item(130, 3).Value = 15.763999999999999
$ws.Cells.Item(130, 4).Value = 150.21

$ws.Cells.Item(131, 1).Value = 20250610
$ws.Cells.Item(131, 2).Value = "A5"
$ws.Cells.Item(131, 3).Value = 15.973000000000001
$ws.Cells.Item(131, 4).Value = 147.834

$ws.Cells.Item(132, 1).Value = 20250610
$ws.Cells.Item(132, 2).Value = "B1"
$ws.Cells.Item(132, 3).Value = 21.608000000000001
$ws.Cells.Item(132, 4).Value = 227.85499999999999

$ws.Cells.Item(133, 1).Value = 20250610
$ws.Cells.Item(133, 2).Value = "B2"
$ws.Cells.Item(133, 3).Value = 20.777000000000001
$ws.Cells.Item(133, 4).Value = 231.61699999999999

$ws.Cells.Item(134, 1).Value = 20250610
$ws.Cells.Item(134, 2).Value = "B3"
$ws.Cells.Item(134, 3).Value = 19.744
$ws.Cells.Item(134, 4).Value = 245.61799999999999

$ws.Cells.Item(135, 1).Value = 20250610
$ws.Cells.Item(135, 2).Value = "B4"
$ws.Cells.Item(135, 3).Value = 21.991
$ws.Cells.Item(135, 4).Value = 279.88400000000001

$ws.Cells.Item(136, 1).Value = 20250610
$ws.Cells.Item(136, 2).Value = "B5"
$ws.Cells.Item(136, 3).Value = 20.222999999999999
$ws.Cells.Item(136, 4).Value = 194.06299999999999

$ws.Cells.Item(137, 1).Value = 20250610
$ws.Cells.Item(137, 2).Value = "C1"
$ws.Cells.Item(137, 3).Value = 19.503
$ws.Cells.Item(137, 4).Value = 264.7

$ws.Cells.Item(138, 1).Value = 20250610
$ws.Cells.Item(138, 2).Value = "C2"
$ws.Cells.Item(138, 3).Value = 20.873000000000001
$ws.Cells.Item(138, 4).Value = 232.108

$ws.Cells.Item(139, 1).Value = 20250610
$ws.Cells.Item(139, 2).Value = "C3"
$ws.Cells.Item(139, 3).Value = 16.806999999999999
$ws.Cells.Item(139, 4).Value = 166.42500000000001

$ws.Cells.Item(140, 1).Value = 20250610
$ws.Cells.Item(140, 2).Value = "C4"
$ws.Cells.Item(140, 3).Value = 18.510999999999999
$ws.Cells.Item(140, 4).Value = 176.79

$ws.Cells.Item(141, 1).Value = 20250610
$ws.Cells.Item(141, 2).Value = "C5"
$ws.Cells.Item(141, 3).Value = 20.67
$ws.Cells.Item(141, 4).Value = 243.31

$ws.Cells.Item(142, 1).Value = 20250610
$ws.Cells.Item(142, 2).Value = "D1"
$ws.Cells.Item(142, 3).Value = 13.461
$ws.Cells.Item(142, 4).Value = 121.627

$ws.Cells.Item(143, 1).Value = 20250610
$ws.Cells.Item(143, 2).Value = "D2"
$ws.Cells.Item(143, 3).Value = 16.422000000000001
$ws.Cells.Item(143, 4).Value = 184.797

$ws.Cells.Item(144, 1).Value = 20250610
$ws.Cells.Item(144, 2).Value = "D3"
$ws.Cells.Item(144, 3).Value = 14.881
$ws.Cells.Item(144, 4).Value = 149.41499999999999

$ws.Cells.Item(145, 1).Value = 20250610
$ws.Cells.Item(145, 2).Value = "D4"
$ws.Cells.Item(145, 3).Value = 18.774000000000001
$ws.Cells.Item(145, 4).Value = 184.85599999999999

$ws.Cells.Item(146, 1).Value = 20250610
$ws.Cells.Item(146, 2).Value = "D5"
$ws.Cells.Item(146, 3).Value = 17.244
$ws.Cells.Item(146, 4).Value = 150.85300000000001

$ws.Cells.Item(147, 1).Value = 20250610
$ws.Cells.Item(147, 2).Value = "E1"
$ws.Cells.Item(147, 3).Value = 25.245999999999999
$ws.Cells.Item(147, 4).Value = 349.45499999999998

$ws.Cells.Item(148, 1).Value = 20250610
$ws.Cells.Item(148, 2).Value = "E2"
$ws.Cells.Item(148, 3).Value = 25.617999999999999
$ws.Cells.Item(148, 4).Value = 291.77800000000002

$ws.Cells.Item(149, 1).Value = 20250610
$ws.Cells.Item(149, 2).Value = "E3"
$ws.Cells.Item(149, 3).Value = 21.58
$ws.Cells.Item(149, 4).Value = 271.65800000000002

$ws.Cells.Item(150, 1).Value = 20250610
$ws.Cells.Item(150, 2).Value = "E4"
$ws.Cells.Item(150, 3).Value = 15.941000000000001
$ws.Cells.Item(150, 4).Value = 183.732

$ws.Cells.Item(151, 1).Value = 20250610
$ws.Cells.Item(151, 2).Value = "E5"
$ws.Cells.Item(151, 3).Value = 22.937000000000001
$ws.Cells.Item(151, 4).Value = 280.86399999999998

# Match the source workbook formatting: explicit black font (style index 1)
# used throughout the data rows.
$ws.Range("A127:D151").Font.Color = 0

# Reproduce the author's final selection/scroll state after entering the data.
$ws.Range("A127:D151").Select() | Out-Null
